$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new purchase row recorded on 2025-05-24 right after the last
# existing data row (row 27).
$newRow = 28

# Column A stores dates as literal text (e.g. "04/17/2025"), not real Excel
# dates, matching the other manually-logged rows in this sheet. Temporarily
# force a Text format so Excel doesn't auto-convert the string into a date
# serial, then clear the formatting again so the new cell doesn't end up
# with a stray style that the original rows don't have.
$ws.Range("A$newRow").NumberFormat = "@"
$ws.Range("A$newRow").Value = "05/24/2025"
$ws.Range("A$newRow").ClearFormats()

$ws.Range("B$newRow").Value = 454.8550000000032
$ws.Range("C$newRow").Value = 0.1099251409789925
$ws.Range("D$newRow").Value = 50
